$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 8).Value = 'kitchens'
$ws.Cells.Item(2, 9).Value = 'distractor'
$ws.Cells.Item(2, 11).Value = 'f'
$ws.Cells.Item(2, 12).Value = 'stimuli/img_ps986.png'
$ws.Cells.Item(2, 13).Value = 90.46428571428571
$ws.Cells.Item(2, 14).Value = 70.60714285714286
$ws.Cells.Item(2, 15).Value = 80.53571428571428
$ws.Cells.Item(2, 16).Value = 28
$ws.Cells.Item(2, 17).Value = 10
$ws.Cells.Item(2, 18).Value = 10
$ws.Cells.Item(2, 19).Value = 10
$ws.Cells.Item(3, 12).Value = 'stimuli/img_z3yzz.png'
$ws.Cells.Item(3, 13).Value = 71.71052631578948
$ws.Cells.Item(3, 14).Value = 49.81578947368421
$ws.Cells.Item(3, 15).Value = 60.76315789473685
$ws.Cells.Item(3, 16).Value = 38
$ws.Cells.Item(3, 17).Value = 5
$ws.Cells.Item(3, 18).Value = 5
$ws.Cells.Item(3, 19).Value = 5
$ws.Cells.Item(4, 12).Value = 'stimuli/img_gbypq.png'
$ws.Cells.Item(4, 13).Value = 76.27500000000001
$ws.Cells.Item(4, 14).Value = 51.925
$ws.Cells.Item(4, 15).Value = 64.09999999999999
$ws.Cells.Item(4, 16).Value = 40
$ws.Cells.Item(4, 17).Value = 6
$ws.Cells.Item(4, 18).Value = 6
$ws.Cells.Item(4, 19).Value = 6
$ws.Cells.Item(5, 12).Value = 'stimuli/img_jivhq.png'
$ws.Cells.Item(5, 13).Value = 37
$ws.Cells.Item(5, 14).Value = 22.26530612244898
$ws.Cells.Item(5, 15).Value = 29.63265306122449
$ws.Cells.Item(5, 16).Value = 49
$ws.Cells.Item(5, 17).Value = 2
$ws.Cells.Item(5, 18).Value = 2
$ws.Cells.Item(5, 19).Value = 2
$ws.Cells.Item(6, 12).Value = 'stimuli/img_anzgh.png'
$ws.Cells.Item(6, 13).Value = 75.10526315789474
$ws.Cells.Item(6, 14).Value = 55.76315789473684
$ws.Cells.Item(6, 15).Value = 65.43421052631579
$ws.Cells.Item(6, 16).Value = 38
$ws.Cells.Item(6, 17).Value = 6
$ws.Cells.Item(6, 18).Value = 6
$ws.Cells.Item(6, 19).Value = 6
$ws.Cells.Item(7, 8).Value = 'bedrooms'
$ws.Cells.Item(7, 9).Value = 'target'
$ws.Cells.Item(7, 11).Value = 'j'
$ws.Cells.Item(7, 12).Value = 'stimuli/img_2pnl2.png'
$ws.Cells.Item(7, 13).Value = 6.621621621621622
$ws.Cells.Item(7, 14).Value = 7.135135135135135
$ws.Cells.Item(7, 15).Value = 6.878378378378379
$ws.Cells.Item(7, 16).Value = 37
$ws.Cells.Item(7, 17).Value = 1
$ws.Cells.Item(7, 18).Value = 1
$ws.Cells.Item(7, 19).Value = 1
$ws.Cells.Item(8, 8).Value = 'bedrooms'
$ws.Cells.Item(8, 9).Value = 'target'
$ws.Cells.Item(8, 11).Value = 'j'
$ws.Cells.Item(8, 12).Value = 'stimuli/img_ic3os.png'
$ws.Cells.Item(8, 13).Value = 84.79069767441861
$ws.Cells.Item(8, 14).Value = 66.16279069767442
$ws.Cells.Item(8, 15).Value = 75.47674418604652
$ws.Cells.Item(8, 16).Value = 43
$ws.Cells.Item(8, 17).Value = 9
$ws.Cells.Item(8, 18).Value = 9
$ws.Cells.Item(8, 19).Value = 9
$ws.Cells.Item(9, 12).Value = 'stimuli/img_xguy9.png'
$ws.Cells.Item(9, 13).Value = 78.21621621621621
$ws.Cells.Item(9, 14).Value = 57.24324324324324
$ws.Cells.Item(9, 15).Value = 67.72972972972973
$ws.Cells.Item(9, 16).Value = 37
$ws.Cells.Item(9, 17).Value = 7
$ws.Cells.Item(9, 18).Value = 7
$ws.Cells.Item(9, 19).Value = 7
$ws.Cells.Item(11, 12).Value = 'stimuli/img_ose78.png'
$ws.Cells.Item(11, 13).Value = 80.19444444444444
$ws.Cells.Item(11, 14).Value = 60.25
$ws.Cells.Item(11, 15).Value = 70.22222222222223
$ws.Cells.Item(11, 16).Value = 36
$ws.Cells.Item(11, 17).Value = 8
$ws.Cells.Item(11, 18).Value = 7
$ws.Cells.Item(11, 19).Value = 7
$ws.Cells.Item(12, 12).Value = 'stimuli/img_cmyvx.png'
$ws.Cells.Item(12, 13).Value = 64.25
$ws.Cells.Item(12, 14).Value = 40.09375
$ws.Cells.Item(12, 15).Value = 52.171875
$ws.Cells.Item(12, 16).Value = 32
$ws.Cells.Item(12, 17).Value = 4
$ws.Cells.Item(12, 18).Value = 4
$ws.Cells.Item(12, 19).Value = 4
$ws.Cells.Item(13, 8).Value = 'living_rooms'
$ws.Cells.Item(13, 9).Value = 'distractor'
$ws.Cells.Item(13, 11).Value = 'f'
$ws.Cells.Item(13, 12).Value = 'stimuli/img_95hiv.png'
$ws.Cells.Item(13, 13).Value = 84.04545454545455
$ws.Cells.Item(13, 14).Value = 67.31818181818181
$ws.Cells.Item(13, 15).Value = 75.68181818181819
$ws.Cells.Item(13, 16).Value = 44
$ws.Cells.Item(13, 17).Value = 9
$ws.Cells.Item(13, 18).Value = 9
$ws.Cells.Item(13, 19).Value = 9
$ws.Cells.Item(14, 12).Value = 'stimuli/img_kzg3h.png'
$ws.Cells.Item(14, 13).Value = 77.02777777777777
$ws.Cells.Item(14, 14).Value = 56.22222222222222
$ws.Cells.Item(14, 15).Value = 66.625
$ws.Cells.Item(14, 16).Value = 36
$ws.Cells.Item(14, 17).Value = 7
$ws.Cells.Item(14, 18).Value = 7
$ws.Cells.Item(14, 19).Value = 7
$ws.Cells.Item(15, 8).Value = 'living_rooms'
$ws.Cells.Item(15, 9).Value = 'distractor'
$ws.Cells.Item(15, 11).Value = 'f'
$ws.Cells.Item(15, 12).Value = 'stimuli/img_zxvl3.png'
$ws.Cells.Item(15, 13).Value = 68.78260869565217
$ws.Cells.Item(15, 14).Value = 47.56521739130435
$ws.Cells.Item(15, 15).Value = 58.17391304347827
$ws.Cells.Item(15, 16).Value = 46
$ws.Cells.Item(16, 8).Value = 'bedrooms'
$ws.Cells.Item(16, 9).Value = 'target'
$ws.Cells.Item(16, 11).Value = 'j'
$ws.Cells.Item(16, 12).Value = 'stimuli/img_72fmj.png'
$ws.Cells.Item(16, 13).Value = 53.87179487179487
$ws.Cells.Item(16, 14).Value = 36.02564102564103
$ws.Cells.Item(16, 15).Value = 44.94871794871795
$ws.Cells.Item(16, 16).Value = 39
$ws.Cells.Item(16, 17).Value = 3
$ws.Cells.Item(16, 18).Value = 3
$ws.Cells.Item(16, 19).Value = 3
$ws.Cells.Item(17, 8).Value = 'living_rooms'
$ws.Cells.Item(17, 9).Value = 'distractor'
$ws.Cells.Item(17, 11).Value = 'f'
$ws.Cells.Item(17, 12).Value = 'stimuli/img_pbsj1.png'
$ws.Cells.Item(17, 13).Value = 73.88636363636364
$ws.Cells.Item(17, 14).Value = 51.52272727272727
$ws.Cells.Item(17, 15).Value = 62.70454545454545
$ws.Cells.Item(17, 16).Value = 44
$ws.Cells.Item(17, 17).Value = 6
$ws.Cells.Item(17, 18).Value = 6
$ws.Cells.Item(17, 19).Value = 6
$ws.Cells.Item(18, 12).Value = 'stimuli/img_cgdyc.png'
$ws.Cells.Item(18, 13).Value = 32.93023255813954
$ws.Cells.Item(18, 14).Value = 14.04651162790698
$ws.Cells.Item(18, 15).Value = 23.48837209302326
$ws.Cells.Item(18, 16).Value = 43
$ws.Cells.Item(18, 17).Value = 1
$ws.Cells.Item(18, 18).Value = 1
$ws.Cells.Item(18, 19).Value = 1
$ws.Cells.Item(19, 12).Value = 'stimuli/img_f4jxo.png'
$ws.Cells.Item(19, 13).Value = 82.91666666666667
$ws.Cells.Item(19, 14).Value = 65.52777777777777
$ws.Cells.Item(19, 15).Value = 74.22222222222223
$ws.Cells.Item(19, 16).Value = 36
$ws.Cells.Item(19, 17).Value = 8
$ws.Cells.Item(19, 18).Value = 8
$ws.Cells.Item(19, 19).Value = 8
$ws.Cells.Item(20, 8).Value = 'bedrooms'
$ws.Cells.Item(20, 9).Value = 'target'
$ws.Cells.Item(20, 11).Value = 'j'
$ws.Cells.Item(20, 12).Value = 'stimuli/img_juob3.png'
$ws.Cells.Item(20, 13).Value = 79.92105263157895
$ws.Cells.Item(20, 14).Value = 59.78947368421053
$ws.Cells.Item(20, 15).Value = 69.85526315789474
$ws.Cells.Item(20, 16).Value = 38
$ws.Cells.Item(20, 17).Value = 7
$ws.Cells.Item(20, 18).Value = 7
$ws.Cells.Item(20, 19).Value = 7
$ws.Cells.Item(21, 8).Value = 'bedrooms'
$ws.Cells.Item(21, 9).Value = 'target'
$ws.Cells.Item(21, 11).Value = 'j'
$ws.Cells.Item(21, 12).Value = 'stimuli/img_1vq1v.png'
$ws.Cells.Item(21, 13).Value = 69.42857142857143
$ws.Cells.Item(21, 14).Value = 46.59523809523809
$ws.Cells.Item(21, 15).Value = 58.01190476190476
$ws.Cells.Item(21, 16).Value = 42
$ws.Cells.Item(21, 17).Value = 5
$ws.Cells.Item(21, 18).Value = 5
$ws.Cells.Item(21, 19).Value = 5
$ws.Cells.Item(22, 12).Value = 'stimuli/img_yteqw.png'
$ws.Cells.Item(22, 13).Value = 66.83783783783784
$ws.Cells.Item(22, 14).Value = 43.78378378378378
$ws.Cells.Item(22, 15).Value = 55.31081081081081
$ws.Cells.Item(22, 16).Value = 37
$ws.Cells.Item(22, 17).Value = 4
$ws.Cells.Item(22, 18).Value = 4
$ws.Cells.Item(22, 19).Value = 4
$ws.Cells.Item(23, 8).Value = 'kitchens'
$ws.Cells.Item(23, 9).Value = 'distractor'
$ws.Cells.Item(23, 11).Value = 'f'
$ws.Cells.Item(23, 12).Value = 'stimuli/img_eppte.png'
$ws.Cells.Item(23, 13).Value = 78.42424242424242
$ws.Cells.Item(23, 14).Value = 57.03030303030303
$ws.Cells.Item(23, 15).Value = 67.72727272727272
$ws.Cells.Item(23, 16).Value = 33
$ws.Cells.Item(23, 17).Value = 7
$ws.Cells.Item(23, 18).Value = 7
$ws.Cells.Item(23, 19).Value = 7
$ws.Cells.Item(24, 12).Value = 'stimuli/img_9pfbj.png'
$ws.Cells.Item(24, 13).Value = 91.27272727272727
$ws.Cells.Item(24, 14).Value = 80.09090909090909
$ws.Cells.Item(24, 15).Value = 85.68181818181819
$ws.Cells.Item(24, 16).Value = 33
$ws.Cells.Item(25, 12).Value = 'stimuli/img_fnu4h.png'
$ws.Cells.Item(25, 13).Value = 85.87179487179488
$ws.Cells.Item(25, 14).Value = 70.71794871794872
$ws.Cells.Item(25, 15).Value = 78.2948717948718
$ws.Cells.Item(25, 16).Value = 39
$ws.Cells.Item(25, 17).Value = 9
$ws.Cells.Item(25, 18).Value = 9
$ws.Cells.Item(25, 19).Value = 9
$ws.Cells.Item(26, 12).Value = 'stimuli/img_t4hvr.png'
$ws.Cells.Item(26, 13).Value = 61.69230769230769
$ws.Cells.Item(26, 14).Value = 39.76923076923077
$ws.Cells.Item(26, 15).Value = 50.73076923076923
$ws.Cells.Item(26, 16).Value = 39
$ws.Cells.Item(26, 17).Value = 3
$ws.Cells.Item(26, 18).Value = 3
$ws.Cells.Item(26, 19).Value = 3
$ws.Cells.Item(27, 12).Value = 'stimuli/img_3bxjb.png'
$ws.Cells.Item(27, 13).Value = 87.28571428571429
$ws.Cells.Item(27, 14).Value = 72.65714285714286
$ws.Cells.Item(27, 15).Value = 79.97142857142858
$ws.Cells.Item(27, 16).Value = 35
$ws.Cells.Item(27, 17).Value = 10
$ws.Cells.Item(27, 18).Value = 10
$ws.Cells.Item(27, 19).Value = 10
